$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macro_taxonomy")
$ws.Select()
